$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 176, shifting the old summary rows (old 176-179) down to 177-180.
$ws.Rows("176:176").Insert()

# Carry the number-format styles used by the other data rows onto the new row.
$ws.Range("D176").NumberFormat = $ws.Range("D175").NumberFormat
$ws.Range("E176").NumberFormat = $ws.Range("E175").NumberFormat
$ws.Range("F176").NumberFormat = $ws.Range("F175").NumberFormat
$ws.Range("G176").NumberFormat = $ws.Range("G175").NumberFormat

# Populate the new data row 176 (same pattern as the preceding data rows).
$ws.Range("A176").Value = 2014
$ws.Range("B176").Value = 8
$ws.Range("C176").Value = 5
$ws.Range("D176").Value = 0.79166666666666663
$ws.Range("E176").Value = 0.875
$ws.Range("F176").Formula = "=(E176-D176)*24*60"
$ws.Range("G176").Formula = "=F176/60"

# Move the sheet's selection to match the author's edit.
$ws.Range("I176").Select()

$wb.Save()
